# Update Observed data output:
# Swap the "conc" and "conc S.D." columns (B and C), including header
# labels and all data values, on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used data row (header in row 1, data starting row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Swap header labels in B1 and C1
$headerB = $ws.Range("B1").Value()
$headerC = $ws.Range("C1").Value()
$ws.Range("B1").Value = $headerC
$ws.Range("C1").Value = $headerB

# Swap the data values in columns B and C for each row
for ($r = 2; $r -le $lastRow; $r++) {
    $valB = $ws.Cells.Item($r, 2).Value()
    $valC = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 2).Value = $valC
    $ws.Cells.Item($r, 3).Value = $valB
}
